$d = $word.ActiveDocument

$d.Content.Find.Execute("960÷7=137, 1", $true, $false, $false, $false, $false, $true, 1, $false, "815÷9=90, 5", 2) | Out-Null
$d.Content.Find.Execute("911÷2=455, 1", $true, $false, $false, $false, $false, $true, 1, $false, "763÷2=381, 1", 2) | Out-Null
$d.Content.Find.Execute("532÷3=177, 1", $true, $false, $false, $false, $false, $true, 1, $false, "305÷5=61, 0", 2) | Out-Null
$d.Content.Find.Execute("880÷3=293, 1", $true, $false, $false, $false, $false, $true, 1, $false, "122÷3=40, 2", 2) | Out-Null
$d.Content.Find.Execute("904÷9=100, 4", $true, $false, $false, $false, $false, $true, 1, $false, "705÷2=352, 1", 2) | Out-Null
$d.Content.Find.Execute("843÷2=421, 1", $true, $false, $false, $false, $false, $true, 1, $false, "651÷3=217, 0", 2) | Out-Null
$d.Content.Find.Execute("679÷5=135, 4", $true, $false, $false, $false, $false, $true, 1, $false, "937÷2=468, 1", 2) | Out-Null
$d.Content.Find.Execute("264÷3=88, 0", $true, $false, $false, $false, $false, $true, 1, $false, "396÷3=132, 0", 2) | Out-Null
$d.Content.Find.Execute("135÷4=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "202÷3=67, 1", 2) | Out-Null
$d.Content.Find.Execute("250÷6=41, 4", $true, $false, $false, $false, $false, $true, 1, $false, "691÷9=76, 7", 2) | Out-Null
$d.Content.Find.Execute("899÷5=179, 4", $true, $false, $false, $false, $false, $true, 1, $false, "960÷2=480, 0", 2) | Out-Null
$d.Content.Find.Execute("343÷6=57, 1", $true, $false, $false, $false, $false, $true, 1, $false, "294÷3=98, 0", 2) | Out-Null
$d.Content.Find.Execute("315÷8=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "453÷5=90, 3", 2) | Out-Null
$d.Content.Find.Execute("909÷2=454, 1", $true, $false, $false, $false, $false, $true, 1, $false, "155÷7=22, 1", 2) | Out-Null
$d.Content.Find.Execute("544÷9=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "512÷9=56, 8", 2) | Out-Null
$d.Content.Find.Execute("643÷9=71, 4", $true, $false, $false, $false, $false, $true, 1, $false, "669÷4=167, 1", 2) | Out-Null
$d.Content.Find.Execute("121÷4=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "394÷9=43, 7", 2) | Out-Null
$d.Content.Find.Execute("861÷2=430, 1", $true, $false, $false, $false, $false, $true, 1, $false, "402÷5=80, 2", 2) | Out-Null
$d.Content.Find.Execute("638÷7=91, 1", $true, $false, $false, $false, $false, $true, 1, $false, "159÷5=31, 4", 2) | Out-Null
$d.Content.Find.Execute("516÷5=103, 1", $true, $false, $false, $false, $false, $true, 1, $false, "113÷7=16, 1", 2) | Out-Null
$d.Content.Find.Execute("513÷3=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "446÷7=63, 5", 2) | Out-Null
$d.Content.Find.Execute("411÷9=45, 6", $true, $false, $false, $false, $false, $true, 1, $false, "460÷2=230, 0", 2) | Out-Null
$d.Content.Find.Execute("349÷5=69, 4", $true, $false, $false, $false, $false, $true, 1, $false, "810÷9=90, 0", 2) | Out-Null
$d.Content.Find.Execute("253÷9=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "354÷5=70, 4", 2) | Out-Null
$d.Content.Find.Execute("320÷7=45, 5", $true, $false, $false, $false, $false, $true, 1, $false, "527÷7=75, 2", 2) | Out-Null
